$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '72.464.45'
$ws.Range('E2').Value = '  +4.51%  '
$ws.Range('D3').Value = '4.056.40'
$ws.Range('E3').Value = '  +4.07%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').Value = '519.20'
$ws.Range('E5').Value = '  -1.84%  '
$ws.Range('D6').Value = '147.28'
$ws.Range('E6').Value = '  +1.93%  '
$ws.Range('D7').Value = '0.734'
$ws.Range('E7').Value = '  +19.97%  '
$ws.Range('D8').Value = '4.051.01'
$ws.Range('E8').Value = '  +4.38%  '
$ws.Range('D9').Value = '0.999'
$ws.Range('D10').Value = '0.772'
$ws.Range('E10').Value = '  +7.63%  '
$ws.Range('E11').Value = '  +1.62%  '
$ws.Range('D12').Value = '0.0000328'
$ws.Range('E12').Value = '  -2.19%  '
$ws.Range('D13').Value = '47.71'
$ws.Range('E13').Value = '  +13.47%  '
$ws.Range('E14').Value = '  +8.82%  '
$ws.Range('D15').Value = '4.702.56'
$ws.Range('E15').Value = '  +4.06%  '
$ws.Range('D16').Value = '4.075.07'
$ws.Range('E16').Value = '  +4.56%  '
$ws.Range('E17').Value = '  +7.09%  '
$ws.Range('D18').Value = '14.11'
$ws.Range('E18').Value = '  +1.00%  '
$ws.Range('E19').Value = '  +0.14%  '
$ws.Range('E20').Value = '  -0.77%  '
$ws.Range('D21').Value = '72.444.97'
$ws.Range('E21').Value = '  +4.60%  '
$ws.Range('D22').Value = '443.44'
$ws.Range('E22').Value = '  +4.11%  '
$ws.Range('D23').Value = '104.69'
$ws.Range('E23').Value = '  +18.81%  '
$ws.Range('E24').Value = '  +5.06%  '
$ws.Range('D25').Value = '14.84'
$ws.Range('E25').Value = '  +5.09%  '
$ws.Range('D26').Value = '4.02'
$ws.Range('E26').Value = '  -0.19%  '
$ws.Range('D27').Value = '11.47'
$ws.Range('E27').Value = '  +0.84%  '
$ws.Range('E28').Value = '  +4.22%  '
$ws.Range('D29').Value = '37.82'
$ws.Range('E29').Value = '  +4.04%  '
$ws.Range('D30').Value = '5.80'
$ws.Range('E30').Value = '  +2.36%  '
$ws.Range('D31').Value = '3.25'
$ws.Range('E31').Value = '  +15.21%  '
$ws.Range('D32').Value = '13.68'
$ws.Range('E32').Value = '  +3.76%  '
$ws.Range('E33').Value = '  +3.56%  '
$ws.Range('D34').Value = '682.50'
$ws.Range('E34').Value = '  -1.07%  '
$ws.Range('D35').Value = '6.82'
$ws.Range('E35').Value = '  +14.85%  '
$ws.Range('D36').Value = '67.04'
$ws.Range('E36').Value = '  -2.73%  '
$ws.Range('D37').Value = '42.43'
$ws.Range('E37').Value = '  +6.36%  '
$ws.Range('D38').Value = '0.0₃0865'
$ws.Range('E38').Value = '  -2.53%  '
$ws.Range('E39').Value = '  -1.39%  '
$ws.Range('E40').Value = '  +6.76%  '
$ws.Range('D41').Value = '0.150'
$ws.Range('E41').Value = '  +0.93%  '
$ws.Range('E42').Value = '  +0.02%  '
$ws.Range('E43').Value = '  +3.77%  '
$ws.Range('D44').Value = '0.999'
$ws.Range('E44').Value = '  -0.10%  '
$ws.Range('D45').Value = '3.27'
$ws.Range('E45').Value = '  +1.59%  '
$ws.Range('D46').Value = '0.159'
$ws.Range('E46').Value = '  +13.49%  '
$ws.Range('E47').Value = '  +2.76%  '
$ws.Range('D48').Value = '2.69'
$ws.Range('E48').Value = '  -4.24%  '
$ws.Range('E49').Value = '  +2.36%  '
$ws.Range('E50').Value = '  +7.70%  '
$ws.Range('E51').Value = '  +2.15%  '
